$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Source:" marker row (uses the newly-introduced shared string) ---
$ws.Cells.Item(204, 2).Value = "Source: 2021-02-30"

# --- Week 10 data (age groups 0-4 .. >=80) ---
$week10 = @(
    @(2021, 10, "0-4",   18750,  5.1),
    @(2021, 10, "5-14",  28009,  8.1),
    @(2021, 10, "15-34", 108796, 7.2),
    @(2021, 10, "35-59", 161342, 7),
    @(2021, 10, "60-79", 89815,  5.7),
    @(2021, 10, ">=80",  44444,  5)
)

$row = 205
foreach ($r in $week10) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# Row 211 is intentionally left blank (matches the existing spacer-row
# pattern between weekly blocks elsewhere in the sheet).

# --- Week 11 data (age groups 0-4 .. >=80) ---
$week11 = @(
    @(2021, 11, "0-4",   28472,  5.1),
    @(2021, 11, "5-14",  41666,  8),
    @(2021, 11, "15-34", 119444, 8.6),
    @(2021, 11, "35-59", 170601, 8.5),
    @(2021, 11, "60-79", 91890,  6.6),
    @(2021, 11, ">=80",  44444,  5.2)
)

$row = 212
foreach ($r in $week11) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# --- Update the frozen-pane scroll position / active selection to the new
#     bottom of the sheet (mirrors the view-state Excel records on save). ---
$win = $excel.ActiveWindow
[void]$ws.Range("E215").Select()
$win.ScrollRow = 194
